# Update countries & provincias Spain
# Applies the refreshed COVID-19 snapshot to the "Pais" sheet:
#  - Peru's case count grew enough to overtake Republica Dominicana,
#    Islandia, Mexico and Panama in the ranking (rows 45-49).
#  - Ruanda's case count grew enough to overtake Kenia and Gibraltar
#    in the ranking (rows 121-123).
#  - A handful of other countries (Estados Unidos, Francia, Austria,
#    Monaco) simply got updated figures without changing rank.
#  - The "updated at" footer timestamp moved from 19:50 to 20:25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $country, $values) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
    $ws.Cells.Item($row, 8).Value = $values[6]
}

# Footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 20:25"

# Top countries whose figures refreshed without a rank change
Set-Row 4  "Estados Unidos" @(205438,16908,8762,192148,4912,475,4528)
Set-Row 9  "Francia"        @(56989,4861,10934,42023,5565,509,4032)
Set-Row 16 "Austria"        @(10663,483,1436,9081,215,18,146)

# Peru overtakes Republica Dominicana, Islandia, Mexico and Panama
Set-Row 45 "Peru"                  @(1323,258,394,899,49,0,30)
Set-Row 46 "Republica Dominicana"  @(1284,175,9,1218,0,6,57)
Set-Row 47 "Islandia"              @(1220,85,225,993,12,0,2)
Set-Row 48 "Mexico"                @(1215,121,35,1151,1,1,29)
Set-Row 49 "Panama"                @(1181,0,9,1142,50,0,29)

# Ruanda overtakes Kenia and Gibraltar
Set-Row 121 "Ruanda"     @(82,7,0,82,0,0,0)
Set-Row 122 "Kenia"      @(81,22,3,77,2,0,1)
Set-Row 123 "Gibraltar"  @(81,12,34,47,0,0,0)

# Monaco: only "Casos criticos" (column F) changes
$ws.Cells.Item(129, 6).Value = 2
